# Build site at 2022-01-09 00:29:46 UTC
# Insert 3 new rows before row 12 (pushing the "Programa resumido:" block and
# everything after it down by 3 rows) and populate them with the new
# "Docentes responsáveis:" info block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows starting at row 12, inheriting formatting from the
# row that used to be there ("Programa resumido:" row, style 1/2/3, no custom height).
$ws.Rows.Item(12).Resize(3).Insert()

# Row 12: label only (column A)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13: first professor, duplicated across B and C like the rest of the sheet
$ws.Range("B13").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C13").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# Row 14: second professor, duplicated across B and C
$ws.Range("B14").Value = "7797767 - Viktor Pastoukhov"
$ws.Range("C14").Value = "7797767 - Viktor Pastoukhov"

# The freshly inserted rows carry empty, styled cells in the columns we did
# not populate (B12/C12, A13, A14) because Insert() duplicates the donor
# row's cell formatting across all three columns. Clear those so the
# untouched cells do not persist as empty styled cells in the sheet.
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
